$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.954.75'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '2.551.66'
$ws.Range("E3").Value = '  +3.06%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -2.23%  '
$ws.Range("D9").Value = '2.549.75'
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("E10").Value = '  -2.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.21%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = '3.006.96'
$ws.Range("E15").Value = '  +3.19%  '
$ws.Range("D16").Value = '62.900.41'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D18").Value = '2.521.08'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '333.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("E21").Value = '  -1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +12.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.55%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  +4.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '410.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.397'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -3.55%  '
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.603'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0518'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("E49").Value = '  +3.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.01%  '
